$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) figures for the two generated
# sheets that mirror this data: "展览" and "全部类型".
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 337
    $ws.Range("F3").Value = 1382
}
